# fix(pipelining): fix testing programs #18
#
# 1. Rename Sheet1 -> "Code"
# 2. Insert a new "Pipelining" worksheet after "Code"
# 3. On "Code": fix instruction @0x000 (MOV R2,#1 ROR 0 -> MOV R2,#0 ROR 0,
#    bits G8/R8 1->0), rewrite the load-immediate syntax for 0x002/0x003
#    ("LOAD Rx, [R0, #n]" -> "LOAD Rx, [R0], #n"), and clear bit Q11 (1->0).
# 4. Populate "Pipelining" with the F/E1/E2/ST pipeline diagram.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename "Sheet1" -> "Code" -------------------------------
$code = $wb.Worksheets.Item(1)
$code.Name = "Code"

# --- Add the new "Pipelining" sheet, right after "Code" ---------------
$pipelining = $wb.Worksheets.Add($null, $code)
$pipelining.Name = "Pipelining"

# === Code sheet fixes ===================================================

# Row 8 (instruction @0x000): MOV R2, #1 ROR 0 -> MOV R2, #0 ROR 0
$code.Range("B8").Value = "MOV R2, #0 ROR 0"
$code.Range("G8").Value = 0
$code.Range("R8").Value = 0

# Row 9 (instruction @0x001): text unchanged (LOAD R1, [R0, #1]!)
$code.Range("B9").Value = "LOAD R1, [R0, #1]!"

# Row 10 (instruction @0x002): syntax fix
$code.Range("B10").Value = "LOAD R3, [R0], #1"

# Row 11 (instruction @0x003): syntax fix + bit fix
$code.Range("B11").Value = "LOAD R2, [R0], #2"
$code.Range("Q11").Value = 0

# Restore selection/active cell on the Code sheet
$code.Range("S11").Select() | Out-Null

# === Pipelining sheet content ===========================================

$pipelining.Columns.Item(2).ColumnWidth = 17.26

# Row 1 - instruction @0x000
$pipelining.Range("A1").Value = "0x000"
$pipelining.Range("B1").Value = "MOV R2, #1 ROR 0"
$pipelining.Range("C1").Value = "F"
$pipelining.Range("D1").Value = "E1"
$pipelining.Range("E1").Value = "E2"

# Row 2 - instruction @0x001
$pipelining.Range("A2").Value = "0x001"
$pipelining.Range("B2").Value = "LOAD R1, [R0, #1]!"
$pipelining.Range("B2").HorizontalAlignment = -4131
$pipelining.Range("D2").Value = "F"
$pipelining.Range("E2").Value = "ST"
$pipelining.Range("F2").Value = "E1"
$pipelining.Range("G2").Value = "E2"

# Row 3 - instruction @0x002
$pipelining.Range("A3").Value = "0x002"
$pipelining.Range("B3").Value = "LOAD R3, [R0], #1"
$pipelining.Range("B3").HorizontalAlignment = -4131
$pipelining.Range("F3").Value = "ST"
$pipelining.Range("G3").Value = "F"
$pipelining.Range("H3").Value = "E1"
$pipelining.Range("I3").Value = "E2"

# Row 4 - instruction @0x003
$pipelining.Range("A4").Value = "0x003"
$pipelining.Range("B4").Value = "LOAD R2, [R0], #2"
$pipelining.Range("B4").HorizontalAlignment = -4131
$pipelining.Range("H4").Value = "ST"
$pipelining.Range("I4").Value = "F"
$pipelining.Range("J4").Value = "E1"
$pipelining.Range("K4").Value = "E2"

# Row 5 - instruction @0x004
$pipelining.Range("A5").Value = "0x004"
$pipelining.Range("B5").Value = "STP"
$pipelining.Range("B5").HorizontalAlignment = -4131
$pipelining.Range("J5").Value = "ST"
$pipelining.Range("K5").Value = "F"
$pipelining.Range("L5").Value = "E1"
$pipelining.Range("M5").Value = "E2"

$pipelining.Range("B27").Select() | Out-Null

# Leave the "Code" sheet as the active/visible tab, matching the source
# workbook (tabSelected="1" on the Code sheet).
$code.Activate() | Out-Null
$code.Range("S11").Select() | Out-Null
